$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footers: both the primary and first-page footer contain the Pearson
# logo picture currently named "image1.png" — rename it to "image2.png".
for ($fi = 1; $fi -le 2; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        for ($k = 1; $k -le $ftr.Range.InlineShapes.Count; $k++) {
            $ishp = $ftr.Range.InlineShapes.Item($k)
            if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp = $ishp.ConvertToShape()
                $shp.Name = "image2.png"
                $shp.ConvertToInlineShape() | Out-Null
            }
        }
    }
}

# Header (first page): the BTEC logo picture is currently named
# "image2.jpg" — rename it to "image1.jpg".
for ($hi = 1; $hi -le 2; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists) {
        for ($k = 1; $k -le $hdr.Range.InlineShapes.Count; $k++) {
            $ishp = $hdr.Range.InlineShapes.Item($k)
            if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp = $ishp.ConvertToShape()
                $shp.Name = "image1.jpg"
                $shp.ConvertToInlineShape() | Out-Null
            }
        }
    }
}
